$p = $ppt.ActivePresentation

# --- Slide 2 ("The Data" -> "The Goal") ---
$s2 = $p.Slides.Item(2)

# Title: "The Data" -> "The Goal"
$titleShp = $s2.Shapes.Item("Titel 1")
$titleShp.TextFrame.TextRange.Text = "The Goal"

# Body placeholder: replace the two "Sourced..."/"Contains..." paragraphs
# with the new 5-paragraph content (three text paragraphs separated by two
# blank paragraphs), each paragraph with no bullet / no indent.
$body = $s2.Shapes.Item("Inhaltsplatzhalter 2")
$tr = $body.TextFrame.TextRange

$para1 = "It is bad when structures collapse."
$para2 = ""
$para3 = "Fibre Reinforced Concrete boasts greater durability, however the makeup of fibres plays a key role."
$para4 = ""
$para5 = "The Aim is to predict the durability of fibre-reinforced concrete based on fibre parameters."

$cr = [char]13
$tr.Text = $para1 + $cr + $para2 + $cr + $para3 + $cr + $para4 + $cr + $para5

# Split paragraph 1 into two runs: "It is bad when " + "structures collapse."
$splitAt = ("It is bad when ").Length + 1
$p1 = $tr.Paragraphs(1, 1)
$runTail = $tr.Characters($splitAt, ("structures collapse.").Length)
$runTail.Font.Size = 17

# Remove bullet / indentation on every paragraph (marL=0 indent=0 buNone
# semantics achieved via turning the bullet off on each paragraph)
for ($i = 1; $i -le 5; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $para.ParagraphFormat.Bullet.Visible = $false
}

# --- Animation timing: click-reveal for paragraphs 1, 3, 5 (0-based: 0,2,4) ---
$seq = $s2.TimeLine.MainSequence
$eff1 = $seq.AddEffect($body, 1, 0, 1)
$eff1.TextRangeStart = 0
$eff1.TextRangeLength = -1
